$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$excel.Goto($ws.Range("AX1"), $true)
